$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - Gross Margin
$ws.Range("D16").Value = 0.5309
$ws.Range("E16").Value = 0.5613
$ws.Range("F16").Value = 0.607
$ws.Range("G16").Value = 0.691

# Row 20 - Free Cash Flow Margin
$ws.Range("B20").Value = 0.3093
$ws.Range("D20").Value = 0.3006
$ws.Range("E20").Value = 0.3483
$ws.Range("F20").Value = 0.3192
$ws.Range("G20").Value = 0.282

# Row 28 - EBITDA Margin
$ws.Range("B28").Value = 0.3796
$ws.Range("D28").Value = 0.3069
$ws.Range("E28").Value = 0.2756
$ws.Range("F28").Value = 0.257
$ws.Range("G28").Value = 0.2571

# Row 29 - Operating Cash Flow Margin
$ws.Range("B29").Value = 0.3266
$ws.Range("D29").Value = 0.3189
$ws.Range("E29").Value = 0.37
$ws.Range("F29").Value = 0.3455
$ws.Range("G29").Value = 0.314
